# Update loading_percent values for the 380 kV case (rows 2-25, columns C-N
# excluding H, J, L, O which remain 0).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newValues = @{
    "C2" = 4.046671948877044
    "D2" = 6.001736749122584
    "E2" = 9.252398853523399
    "F2" = 33.21051778171408
    "G2" = 3.661432876531466
    "I2" = 26.74732431446872
    "K2" = 16.44031701738331
    "M2" = 16.87421503032314
    "N2" = 18.59966736151639
    "C3" = 4.04135363498933
    "D3" = 6.030787101149486
    "E3" = 9.15387381791747
    "F3" = 32.97055671414231
    "G3" = 3.665237526062143
    "I3" = 26.6256824473975
    "K3" = 15.9831664879608
    "M3" = 16.60665352026527
    "N3" = 18.67559853304753
    "C4" = 4.038046313733309
    "D4" = 6.049371864444581
    "E4" = 9.096084747368538
    "F4" = 32.83382516377936
    "G4" = 3.667692909754605
    "I4" = 26.55928996400788
    "K4" = 15.70056159101154
    "M4" = 16.44527072832233
    "N4" = 18.72420695418237
    "C5" = 4.036688408236035
    "D5" = 6.057133949612965
    "E5" = 9.073238684846725
    "F5" = 32.78081675403172
    "G5" = 3.668723619568194
    "I5" = 26.53433471359065
    "K5" = 15.58510944251621
    "M5" = 16.3803217913605
    "N5" = 18.7445163812341
    "C6" = 4.036462327520802
    "D6" = 6.05843424980782
    "E6" = 9.069488240168694
    "F6" = 32.77217964625891
    "G6" = 3.668896590865731
    "I6" = 26.53031810721129
    "K6" = 15.56592709282868
    "M6" = 16.36958884280502
    "N6" = 18.74791905868233
    "C7" = 4.03802804114252
    "D7" = 6.04947578204706
    "E7" = 9.095773759932973
    "F7" = 32.83309924427751
    "G7" = 3.667706688153761
    "I7" = 26.55894488911507
    "K7" = 15.69900546929927
    "M7" = 16.44439138778987
    "N7" = 18.72447882322642
    "C8" = 4.044847106049335
    "D8" = 6.011598493948291
    "E8" = 9.217880227984745
    "F8" = 33.12560613017658
    "G8" = 3.662720027748371
    "I8" = 26.70366981553876
    "K8" = 16.28320105559553
    "M8" = 16.78140987718592
    "N8" = 18.6254371327962
    "C9" = 4.057871184766486
    "D9" = 5.943225292350762
    "E9" = 9.477704745503909
    "F9" = 33.78116257855807
    "G9" = 3.653882481195824
    "I9" = 27.05252557242729
    "K9" = 17.4056817054261
    "M9" = 17.4614842932893
    "N9" = 18.44690233063252
    "C10" = 4.067207306095294
    "D10" = 5.896549648090618
    "E10" = 9.679457644065501
    "F10" = 34.30946413729905
    "G10" = 3.647955771991425
    "I10" = 27.34723173547359
    "K10" = 18.20642623431412
    "M10" = 17.96769623729232
    "N10" = 18.3251871522941
    "C11" = 4.071399826126596
    "D11" = 5.876079599738928
    "E11" = 9.773245513515535
    "F11" = 34.55914096171638
    "G11" = 3.645380886349955
    "I11" = 27.48930871000196
    "K11" = 18.56366913872255
    "M11" = 18.19834289136634
    "N11" = 18.2718455639287
    "C12" = 4.072979225068574
    "D12" = 5.868437205198366
    "E12" = 9.809019025399577
    "F12" = 34.65495947129221
    "G12" = 3.644423147971407
    "I12" = 27.54422970417544
    "K12" = 18.6978022577826
    "M12" = 18.28564954261261
    "N12" = 18.25193631472638
    "C13" = 4.072639444638853
    "D13" = 5.870078284787591
    "E13" = 9.801303575503024
    "F13" = 34.63426789382102
    "G13" = 3.644628645793425
    "I13" = 27.53235224337061
    "K13" = 18.66896736193799
    "M13" = 18.26684937241904
    "N13" = 18.25621125073764
    "C14" = 4.071529931052152
    "D14" = 5.875448670502453
    "E14" = 9.776183627273802
    "F14" = 34.56699895836219
    "G14" = 3.645301746272565
    "I14" = 27.49380485320254
    "K14" = 18.57472808254825
    "M14" = 18.20552683711335
    "N14" = 18.27020181326996
    "C15" = 4.070849240728087
    "D15" = 5.878752387911566
    "E15" = 9.760829622653516
    "F15" = 34.52595813510035
    "G15" = 3.645716291286007
    "I15" = 27.47033824192682
    "K15" = 18.51685055847223
    "M15" = 18.16795801784196
    "N15" = 18.27880916697621
    "C16" = 4.066932185493846
    "D16" = 5.897902719641841
    "E16" = 9.673366189019246
    "F16" = 34.29332912366694
    "G16" = 3.648126476000779
    "I16" = 27.33810517853725
    "K16" = 18.18292626652116
    "M16" = 17.95262282859245
    "N16" = 18.32871381861643
    "C17" = 4.06451496327286
    "D17" = 5.909845845924115
    "E17" = 9.620203011007753
    "F17" = 34.15296236258135
    "G17" = 3.649636008266167
    "I17" = 27.25901443108546
    "K17" = 17.9761747049515
    "M17" = 17.82055307098619
    "N17" = 18.35984688849613
    "C18" = 4.063119530632344
    "D18" = 5.916787059154148
    "E18" = 9.589815973074163
    "F18" = 34.07311174951727
    "G18" = 3.65051566640585
    "I18" = 27.21427994448227
    "K18" = 17.85660335765821
    "M18" = 17.74462975926527
    "N18" = 18.3779447085864
    "C19" = 4.062646194965591
    "D19" = 5.919149590401189
    "E19" = 9.579561187063238
    "F19" = 34.04622982139698
    "G19" = 3.650815467571738
    "I19" = 27.19926448977161
    "K19" = 17.81601074428741
    "M19" = 17.71893295533483
    "N19" = 18.38410515186437
    "C20" = 4.064772812352407
    "D20" = 5.908567048233399
    "E20" = 9.625842775414165
    "F20" = 34.16781360831547
    "G20" = 3.649474135386898
    "I20" = 27.26735571936639
    "K20" = 17.99825241399721
    "M20" = 17.83460865205402
    "N20" = 18.35651297311687
    "C21" = 4.071856048072465
    "D21" = 5.873868299489305
    "E21" = 9.78355520803979
    "F21" = 34.5867235694401
    "G21" = 3.645103571331044
    "I21" = 27.50509704648542
    "K21" = 18.60244060727741
    "M21" = 18.22354036151663
    "N21" = 18.26608458556701
    "C22" = 4.076437284478289
    "D22" = 5.851826735786432
    "E22" = 9.888119250873073
    "F22" = 34.8678798231925
    "G22" = 3.642348027012046
    "I22" = 27.66698326584378
    "K22" = 18.99057178763337
    "M22" = 18.47749192279025
    "N22" = 18.20867438622344
    "C23" = 4.073996713310073
    "D23" = 5.8635327081052
    "E23" = 9.832185384676224
    "F23" = 34.71717156234465
    "G23" = 3.643809520691827
    "I23" = 27.57999749475104
    "K23" = 18.78407787890014
    "M23" = 18.34200297684669
    "N23" = 18.23916113688768
    "C24" = 4.064656256631864
    "D24" = 5.909144959173005
    "E24" = 9.623292483706981
    "F24" = 34.1610967129378
    "G24" = 3.64954728132815
    "I24" = 27.26358232850524
    "K24" = 17.98827327368555
    "M24" = 17.82825410053871
    "N24" = 18.35801961629755
    "C25" = 4.0543872278063
    "D25" = 5.961094271185639
    "E25" = 9.405390827466489
    "F25" = 33.59536551747615
    "G25" = 3.656173287463867
    "I25" = 26.95130384816489
    "K25" = 17.10555978461937
    "M25" = 17.27598558566489
    "N25" = 18.4935326854059
}

foreach ($cellRef in $newValues.Keys) {
    $ws.Range($cellRef).Value = $newValues[$cellRef]
}

Write-Host "Updated" $newValues.Keys.Count "cells with new loading_percent values"
